$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range's last row
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column E holds the "Category" values. Rename "Career" category to "Income"
# for all rows (secondary income sources previously categorized as "Career").
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value2 -eq "Career") {
        $cell.Value2 = "Income"
    }
}
